$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # _set_RESOURCES
$ws2 = $wb.Worksheets.Item(2)   # _set_PRODUCTS
$ws3 = $wb.Worksheets.Item(3)   # _set_PRODUCT_DATA

# ----- Sheet1: _set_RESOURCES -----
# Rename header from "r_Names" to "resources_Name"; the rest of the column
# (low energy / avg energy / high energy) keeps its existing values.
$ws1.Range("A1").Value = "resources_Name"

# ----- Sheet2: _set_PRODUCTS -----
# Rename header from "p_Names" to "products_Name" and add a new
# "products_Aggregation" header in B1, copying A1's header formatting.
$ws2.Range("A1").Value = "products_Name"
$ws2.Range("B1").Value = "products_Aggregation"
$ws2.Range("A1").Copy()
$ws2.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ----- Sheet3: _set_PRODUCT_DATA -----
# Rename headers and add a new "product_data_Aggregation" header in C1,
# copying A1's header formatting.
$ws3.Range("A1").Value = "product_data_Name"
$ws3.Range("B1").Value = "product_data_category"
$ws3.Range("C1").Value = "product_data_Aggregation"
$ws3.Range("A1").Copy()
$ws3.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ----- Selections -----
# Apply sheet3's selection last so it remains the active/tab-selected sheet,
# matching the workbook's original "active tab" state.
$ws1.Range("A2:A4").Select()
$ws3.Range("A2:B4").Select()
